# #193 tambah tab untuk setiap ujian
# Add one tab per exam: IPU, IBK, IP, PPKP (IBK reuses the original "IPU"
# sheet's content/sheetId, a fresh "IPU" sheet is inserted before it, and
# "IP" / "PPKP" are appended as further copies). The "IBK" tab ends up
# selected.

$wb = $excel.ActiveWorkbook

# The workbook starts with a single sheet named "IPU". Repurpose it as the
# "IBK" sheet (keeps its original sheetId / internal identity).
$ipu = $wb.Worksheets.Item("IPU")
$ipu.Name = "IBK"

# Insert a fresh copy of that sheet immediately before "IBK" and name it
# "IPU" -- this becomes the new first tab.
$ibk = $wb.Worksheets.Item("IBK")
$ibk.Copy($ibk, $null)
$wb.Worksheets.Item(1).Name = "IPU"

# Append a copy after "IBK" named "IP".
$ibk = $wb.Worksheets.Item("IBK")
$ibk.Copy($null, $ibk)
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "IP"

# Append a copy after "IP" named "PPKP".
$ip = $wb.Worksheets.Item("IP")
$ip.Copy($null, $ip)
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "PPKP"

# "IBK" is the selected/active tab.
$wb.Worksheets.Item("IBK").Activate()
